$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

$ws = $wb.Worksheets.Item("ALERTS")
Set-TextCell $ws "A2" "2026-01-30"
Set-TextCell $ws "B2" "15:56:30"
Set-TextCell $ws "C2" "15:00"
Set-TextCell $ws "D2" "Living Room"
Set-TextCell $ws "E2" "CRITICAL"
Set-TextCell $ws "F2" "FALL_DETECTED"

$ws = $wb.Worksheets.Item("mmWave")
Set-TextCell $ws "A8" "2026-01-30"
Set-TextCell $ws "B8" "15:55:59"
Set-TextCell $ws "C8" "15:00"
Set-TextCell $ws "D8" "Living Room"
Set-TextCell $ws "E8" "PRESENCE_DETECTED"
Set-TextCell $ws "F8" "Active"
Set-TextCell $ws "A9" "2026-01-30"
Set-TextCell $ws "B9" "15:56:05"
Set-TextCell $ws "C9" "15:00"
Set-TextCell $ws "D9" "Living Room"
Set-TextCell $ws "E9" "PRESENCE_DETECTED"
Set-TextCell $ws "F9" "Active"
Set-TextCell $ws "A10" "2026-01-30"
Set-TextCell $ws "B10" "15:56:31"
Set-TextCell $ws "C10" "15:00"
Set-TextCell $ws "D10" "Living Room"
Set-TextCell $ws "E10" "PRESENCE_DETECTED"
Set-TextCell $ws "F10" "Active"
Set-TextCell $ws "A11" "2026-01-30"
Set-TextCell $ws "B11" "15:56:47"
Set-TextCell $ws "C11" "15:00"
Set-TextCell $ws "D11" "Living Room"
Set-TextCell $ws "E11" "PRESENCE_DETECTED"
Set-TextCell $ws "F11" "Active"
Set-TextCell $ws "A12" "2026-01-30"
Set-TextCell $ws "B12" "15:56:57"
Set-TextCell $ws "C12" "15:00"
Set-TextCell $ws "D12" "Living Room"
Set-TextCell $ws "E12" "PRESENCE_DETECTED"
Set-TextCell $ws "F12" "Active"

$ws = $wb.Worksheets.Item("PIR")
Set-TextCell $ws "A19" "2026-01-30"
Set-TextCell $ws "B19" "15:55:59"
Set-TextCell $ws "C19" "15:00"
Set-TextCell $ws "D19" "Bathroom"
Set-TextCell $ws "E19" "No Motion"
Set-TextCell $ws "F19" "Inactive"
Set-TextCell $ws "A20" "2026-01-30"
Set-TextCell $ws "B20" "15:56:02"
Set-TextCell $ws "C20" "15:00"
Set-TextCell $ws "D20" "Bathroom"
Set-TextCell $ws "E20" "No Motion"
Set-TextCell $ws "F20" "Inactive"
Set-TextCell $ws "A21" "2026-01-30"
Set-TextCell $ws "B21" "15:56:07"
Set-TextCell $ws "C21" "15:00"
Set-TextCell $ws "D21" "Bathroom"
Set-TextCell $ws "E21" "No Motion"
Set-TextCell $ws "F21" "Inactive"
Set-TextCell $ws "A22" "2026-01-30"
Set-TextCell $ws "B22" "15:56:12"
Set-TextCell $ws "C22" "15:00"
Set-TextCell $ws "D22" "Bathroom"
Set-TextCell $ws "E22" "No Motion"
Set-TextCell $ws "F22" "Inactive"
Set-TextCell $ws "A23" "2026-01-30"
Set-TextCell $ws "B23" "15:56:17"
Set-TextCell $ws "C23" "15:00"
Set-TextCell $ws "D23" "Bathroom"
Set-TextCell $ws "E23" "No Motion"
Set-TextCell $ws "F23" "Inactive"
Set-TextCell $ws "A24" "2026-01-30"
Set-TextCell $ws "B24" "15:56:22"
Set-TextCell $ws "C24" "15:00"
Set-TextCell $ws "D24" "Bathroom"
Set-TextCell $ws "E24" "No Motion"
Set-TextCell $ws "F24" "Inactive"
Set-TextCell $ws "A25" "2026-01-30"
Set-TextCell $ws "B25" "15:56:27"
Set-TextCell $ws "C25" "15:00"
Set-TextCell $ws "D25" "Bathroom"
Set-TextCell $ws "E25" "No Motion"
Set-TextCell $ws "F25" "Inactive"
Set-TextCell $ws "A26" "2026-01-30"
Set-TextCell $ws "B26" "15:56:31"
Set-TextCell $ws "C26" "15:00"
Set-TextCell $ws "D26" "Living Room"
Set-TextCell $ws "E26" "RECOVERY_DETECTION"
Set-TextCell $ws "F26" "Inactive"
Set-TextCell $ws "A27" "2026-01-30"
Set-TextCell $ws "B27" "15:56:32"
Set-TextCell $ws "C27" "15:00"
Set-TextCell $ws "D27" "Bathroom"
Set-TextCell $ws "E27" "No Motion"
Set-TextCell $ws "F27" "Inactive"
Set-TextCell $ws "A28" "2026-01-30"
Set-TextCell $ws "B28" "15:56:37"
Set-TextCell $ws "C28" "15:00"
Set-TextCell $ws "D28" "Bathroom"
Set-TextCell $ws "E28" "No Motion"
Set-TextCell $ws "F28" "Inactive"
Set-TextCell $ws "A29" "2026-01-30"
Set-TextCell $ws "B29" "15:56:42"
Set-TextCell $ws "C29" "15:00"
Set-TextCell $ws "D29" "Bathroom"
Set-TextCell $ws "E29" "No Motion"
Set-TextCell $ws "F29" "Inactive"
Set-TextCell $ws "A30" "2026-01-30"
Set-TextCell $ws "B30" "15:56:48"
Set-TextCell $ws "C30" "15:00"
Set-TextCell $ws "D30" "Bathroom"
Set-TextCell $ws "E30" "No Motion"
Set-TextCell $ws "F30" "Inactive"
Set-TextCell $ws "A31" "2026-01-30"
Set-TextCell $ws "B31" "15:56:53"
Set-TextCell $ws "C31" "15:00"
Set-TextCell $ws "D31" "Bathroom"
Set-TextCell $ws "E31" "No Motion"
Set-TextCell $ws "F31" "Inactive"
Set-TextCell $ws "A32" "2026-01-30"
Set-TextCell $ws "B32" "15:56:58"
Set-TextCell $ws "C32" "15:00"
Set-TextCell $ws "D32" "Bathroom"
Set-TextCell $ws "E32" "No Motion"
Set-TextCell $ws "F32" "Inactive"

$ws = $wb.Worksheets.Item("Humidity")
Set-TextCell $ws "A17" "2026-01-30"
Set-TextCell $ws "B17" "15:56:00"
Set-TextCell $ws "C17" "15:00"
Set-TextCell $ws "D17" "Bathroom"
Set-TextCell $ws "E17" "87.8%"
Set-TextCell $ws "F17" "Active"
Set-TextCell $ws "A18" "2026-01-30"
Set-TextCell $ws "B18" "15:56:03"
Set-TextCell $ws "C18" "15:00"
Set-TextCell $ws "D18" "Bathroom"
Set-TextCell $ws "E18" "87.8%"
Set-TextCell $ws "F18" "Active"
Set-TextCell $ws "A19" "2026-01-30"
Set-TextCell $ws "B19" "15:56:08"
Set-TextCell $ws "C19" "15:00"
Set-TextCell $ws "D19" "Bathroom"
Set-TextCell $ws "E19" "87.8%"
Set-TextCell $ws "F19" "Active"
Set-TextCell $ws "A20" "2026-01-30"
Set-TextCell $ws "B20" "15:56:18"
Set-TextCell $ws "C20" "15:00"
Set-TextCell $ws "D20" "Bathroom"
Set-TextCell $ws "E20" "87.8%"
Set-TextCell $ws "F20" "Active"
Set-TextCell $ws "A21" "2026-01-30"
Set-TextCell $ws "B21" "15:56:23"
Set-TextCell $ws "C21" "15:00"
Set-TextCell $ws "D21" "Bathroom"
Set-TextCell $ws "E21" "86.9%"
Set-TextCell $ws "F21" "Active"
Set-TextCell $ws "A22" "2026-01-30"
Set-TextCell $ws "B22" "15:56:30"
Set-TextCell $ws "C22" "15:00"
Set-TextCell $ws "D22" "Bathroom"
Set-TextCell $ws "E22" "87.7%"
Set-TextCell $ws "F22" "Active"
Set-TextCell $ws "A23" "2026-01-30"
Set-TextCell $ws "B23" "15:56:33"
Set-TextCell $ws "C23" "15:00"
Set-TextCell $ws "D23" "Bathroom"
Set-TextCell $ws "E23" "86.8%"
Set-TextCell $ws "F23" "Active"
Set-TextCell $ws "A24" "2026-01-30"
Set-TextCell $ws "B24" "15:56:38"
Set-TextCell $ws "C24" "15:00"
Set-TextCell $ws "D24" "Bathroom"
Set-TextCell $ws "E24" "87.8%"
Set-TextCell $ws "F24" "Active"
Set-TextCell $ws "A25" "2026-01-30"
Set-TextCell $ws "B25" "15:56:43"
Set-TextCell $ws "C25" "15:00"
Set-TextCell $ws "D25" "Bathroom"
Set-TextCell $ws "E25" "86.8%"
Set-TextCell $ws "F25" "Active"
Set-TextCell $ws "A26" "2026-01-30"
Set-TextCell $ws "B26" "15:56:48"
Set-TextCell $ws "C26" "15:00"
Set-TextCell $ws "D26" "Bathroom"
Set-TextCell $ws "E26" "87.7%"
Set-TextCell $ws "F26" "Active"
Set-TextCell $ws "A27" "2026-01-30"
Set-TextCell $ws "B27" "15:56:58"
Set-TextCell $ws "C27" "15:00"
Set-TextCell $ws "D27" "Bathroom"
Set-TextCell $ws "E27" "87.7%"
Set-TextCell $ws "F27" "Active"
